# Append 4 new order rows (21-24) to the "Orders" sheet for the new
# customer "JohnDoe" (orders for "Sakamoto days !" and "One Piece"),
# mirroring the rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Column layout: A User Name | B City | C Phone | D Order Status |
#                E Placed At | F Manga Title | G Volume Number |
#                H Quantity  | I Unit Price  | J Item Total | K Order Total
$rows = @(
    @{ Row = 21; A = "JohnDoe"; B = "test";   C = "07893830932"; D = "PENDING"; E = "2025-08-02T12:03:50.585Z"; F = "Sakamoto days !"; G = 3; H = 1;  I = 108.9; J = 108.9;  K = 326.7 },
    @{ Row = 22; A = "JohnDoe"; B = "test";   C = "07893830932"; D = "PENDING"; E = "2025-08-02T12:03:50.585Z"; F = "Sakamoto days !"; G = 6; H = 2;  I = 108.9; J = 217.8;  K = 326.7 },
    @{ Row = 23; A = "JohnDoe"; B = "tanger"; C = "93932979430"; D = "PENDING"; E = "2025-08-02T12:04:59.577Z"; F = "One Piece";        G = 4; H = 42; I = 108.9; J = 4573.8; K = 9365.4 },
    @{ Row = 24; A = "JohnDoe"; B = "tanger"; C = "93932979430"; D = "PENDING"; E = "2025-08-02T12:04:59.577Z"; F = "Sakamoto days !"; G = 2; H = 44; I = 108.9; J = 4791.6; K = 9365.4 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B

    # The phone number is made of digits only (and can carry a leading
    # zero), so format the cell as text first to stop Excel from
    # reinterpreting it as a number and dropping the leading zero.
    $cPhone = $ws.Cells.Item($row, 3)
    $cPhone.NumberFormat = "@"
    $cPhone.Value = $r.C

    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
}
